# Insert a new data row at row 161 (pushing the existing rows 161-248 down
# to 162-249) and populate it with the new "Arveja Verde" price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("161:161").Insert()

$ws.Range("A161").Value = 6
$ws.Range("B161").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C161").Value = "Metropolitana"
$ws.Range("D161").Value = 44806
$ws.Range("E161").Value = 13
$ws.Range("F161").Value = 100112022
$ws.Range("G161").Value = "Arveja Verde"
$ws.Range("H161").Value = "Perfection"
$ws.Range("I161").Value = "Primera"
$ws.Range("J161").Value = 260
$ws.Range("K161").Value = 36000
$ws.Range("L161").Value = 38000
$ws.Range("M161").Value = 36769
$ws.Range("N161").Value = "`$/malla 25 kilos"
$ws.Range("O161").Value = "Provincia de Huasco"
$ws.Range("P161").Value = 1471
$ws.Range("Q161").Value = 25
$ws.Range("R161").Value = "Hortaliza"
